$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 value (22.1 -> 25.1)
$ws.Range("C2").Value = 25.1

# Update A11 and A30: "ליאם דיין " -> "ליאם דיין" (new shared string)
$ws.Range("A11").Value = "ליאם דיין"
$ws.Range("A30").Value = "ליאם דיין"

# Append new ranking rows 236-296
$names = @(
    "אביב ואסקז",
    "ליאם דיין",
    "יובל סטרוזר",
    "יהלי גודר",
    "תאיו ורד",
    "שלו דיין",
    "שלו דיין",
    "יובל סטרוזר",
    "הגר אגמון",
    "דן פימה",
    "דפנה ברגשטיין",
    "תומר ששון",
    "יהלי דוייב",
    "אן מרש",
    "איתי בסטקר",
    "אורי שטרנברג",
    "יולי יערי תליו",
    "ירון גלפנד",
    "הילה שולויס",
    "דפנה ברגשטיין",
    "יהלי דוייב",
    "רומי הרשקוביץ",
    "עדן ורד מרי",
    "איתי בסטקר",
    "הילה שולויס",
    "איתי הראל",
    "ליהי בראל",
    "יובל סטרוזר",
    "יולי יערי תליו",
    "ליאם דיין",
    "ליאם דיין",
    "איתי בסטקר",
    "דן פימה",
    "הגר אגמון",
    "איתי הראל",
    "דפנה ברגשטיין",
    "יהלי דוייב",
    "יובל סטרוזר",
    "אורי שטרנברג",
    "מעיין סטרוזר",
    "שלו דיין",
    "ירון גלפנד",
    "ירון גלפנד",
    "איתי הראל",
    "הגר אגמון",
    "תומר ששון",
    "אורי שטרנברג",
    "שלו דיין",
    "ירון גלפנד",
    "שלו דיין",
    "שלו דיין",
    "עדן ורד מרי",
    "רומי הרשקוביץ",
    "תאיו ורד",
    "הילה שולויס",
    "ליאם דיין",
    "לינוי קוסטיקה",
    "איתי הראל",
    "ליהי בראל",
    "תאיו ורד",
    "ליהי בראל"
)
$points = @(
    1,
    1,
    1,
    1,
    1,
    1,
    6,
    6,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    6,
    6,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    6,
    6,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    6,
    6,
    1,
    1,
    1,
    1,
    1,
    6,
    6,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    6,
    6
)

$startRow = 236
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $points[$i]
}

# Update selection to C4
[void]$ws.Range("C4").Select()
